# Insert a new weekly price record at row 52 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 52..113 down to 53..114.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44789
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = 100112035
$ws.Range("G52").Value = "Bruselas (repollito)"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 45
$ws.Range("K52").Value = 24000
$ws.Range("L52").Value = 24000
$ws.Range("M52").Value = 24000
$ws.Range("N52").Value = "$/malla 10 kilos"
$ws.Range("O52").Value = "Provincia de Quillota"
$ws.Range("P52").Value = 2400
$ws.Range("Q52").Value = 10
$ws.Range("R52").Value = "Hortaliza"
